$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "75.867.08"
$ws.Range("E2").Value = "  +2.26%  "
$ws.Range("D3").Value = "2.846.76"
$ws.Range("E3").Value = "  +7.07%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "194.12"
$ws.Range("E5").Value = "  +4.71%  "
$ws.Range("D6").Value = "599.27"
$ws.Range("E6").Value = "  +2.81%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "0.551"
$ws.Range("E8").Value = "  +3.81%  "
$ws.Range("D9").Value = "0.194"
$ws.Range("E9").Value = "  +1.00%  "
$ws.Range("D10").Value = "2.842.48"
$ws.Range("D11").Value = "0.391"
$ws.Range("E11").Value = "  +10.54%  "
$ws.Range("E12").Value = "  -2.10%  "
$ws.Range("D13").Value = "4.93"
$ws.Range("E13").Value = "  +5.00%  "
$ws.Range("D14").Value = "3.370.63"
$ws.Range("E14").Value = "  +7.06%  "
$ws.Range("D15").Value = "75.808.74"
$ws.Range("E15").Value = "  +2.25%  "
$ws.Range("D16").Value = "27.61"
$ws.Range("E16").Value = "  +5.27%  "
$ws.Range("D17").Value = "0.0000189"
$ws.Range("E17").Value = "  +2.14%  "
$ws.Range("D18").Value = "2.843.06"
$ws.Range("E18").Value = "  +6.23%  "
$ws.Range("D19").Value = "9.13"
$ws.Range("E19").Value = "  -1.95%  "
$ws.Range("D20").Value = "12.47"
$ws.Range("E20").Value = "  +4.88%  "
$ws.Range("D21").Value = "384.98"
$ws.Range("E21").Value = "  +4.20%  "
$ws.Range("D22").Value = "2.32"
$ws.Range("E22").Value = "  +3.08%  "
$ws.Range("D23").Value = "4.15"
$ws.Range("E23").Value = "  +2.49%  "
$ws.Range("D24").Value = "71.54"
$ws.Range("E24").Value = "  +3.10%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").Value = "4.24"
$ws.Range("E26").Value = "  +3.60%  "
$ws.Range("D27").Value = "2.986.35"
$ws.Range("E27").Value = "  +6.72%  "
$ws.Range("D28").Value = "9.76"
$ws.Range("E28").Value = "  +4.92%  "
$ws.Range("D29").Value = "0.0000105"
$ws.Range("E29").Value = "  +12.15%  "
$ws.Range("D30").Value = "0.997"
$ws.Range("E30").Value = "  -0.33%  "
$ws.Range("E31").Value = "  +3.26%  "
$ws.Range("D32").Value = "521.74"
$ws.Range("E32").Value = "  +1.18%  "
$ws.Range("D33").Value = "7.75"
$ws.Range("E33").Value = "  +1.23%  "
$ws.Range("E34").Value = "  +5.07%  "
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("D36").Value = "166.10"
$ws.Range("E36").Value = "  +1.66%  "
$ws.Range("D37").Value = "20.00"
$ws.Range("E37").Value = "  +4.89%  "
$ws.Range("D38").Value = "0.119"
$ws.Range("E38").Value = "  +1.08%  "
$ws.Range("E39").Value = "  +0.45%  "
$ws.Range("D40").Value = "186.46"
$ws.Range("E40").Value = "  +9.73%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").Value = "5.13"
$ws.Range("E42").Value = "  +4.23%  "
$ws.Range("D43").Value = "0.345"
$ws.Range("E43").Value = "  +6.02%  "
$ws.Range("E44").Value = "  +2.66%  "
$ws.Range("D45").Value = "1.24"
$ws.Range("E45").Value = "  +4.75%  "
$ws.Range("D46").Value = "40.11"
$ws.Range("E46").Value = "  +2.76%  "
$ws.Range("D47").Value = "0.0886"
$ws.Range("E47").Value = "  +5.38%  "
$ws.Range("D48").Value = "2.39"
$ws.Range("E48").Value = "  +2.60%  "
$ws.Range("D49").Value = "0.576"
$ws.Range("E49").Value = "  +9.99%  "
$ws.Range("D50").Value = "3.77"
$ws.Range("E50").Value = "  +4.30%  "
$ws.Range("D51").Value = "0.656"
$ws.Range("E51").Value = "  +11.04%  "
